# Automatic update of files.
#
# The underlying data feed for this sheet was refreshed:
#   * Column C ("Förändrad") is bumped from 46070 to 46072 for every data
#     row (rows 2-23).
#   * A handful of rows were re-sorted by the source feed, which shows up
#     as a left-rotation of the row content within a few small groups of
#     consecutive rows (5-6, 7-8, 14-16, 17-20). Every column (including
#     the Markägare cell and the six HYPERLINK formulas) moves together
#     with its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    $ws.Range($addr).Value = $text
}

function Set-Num($addr, $num) {
    $ws.Range($addr).Value = $num
}

function Set-Link($addr, $target, $label) {
    $formula = '=HYPERLINK("https://klasma.github.io/Logging_1286/' + $target + '", "' + $label + '")'
    $ws.Range($addr).Formula = $formula
}

function Clear-Cell($addr) {
    $ws.Range($addr).ClearContents() | Out-Null
}

# ---------------------------------------------------------------------
# 1) Bump "Förändrad" (column C) to 46072 for every data row (2-23)
# ---------------------------------------------------------------------
for ($r = 2; $r -le 23; $r++) {
    Set-Num ("C" + $r) 46072
}

# ---------------------------------------------------------------------
# 2) Rewrite the rows whose content was re-ordered by the source feed.
#    Each entry below is the FULL new content of that row (taken from the
#    row that rotated into its place).
# ---------------------------------------------------------------------

# --- Row 5 (becomes old row 6: "A 61558-2023") ---
Set-Text "A5" "A 61558-2023"
Set-Num  "B5" 45265
Set-Text "F5" "Övriga statliga verk och myndigheter"
Set-Num  "G5" 1.5
Set-Num  "H5" 0
Set-Num  "I5" 0
Set-Num  "J5" 0
Set-Num  "K5" 0
Set-Num  "L5" 1
Set-Num  "M5" 0
Set-Num  "N5" 0
Set-Num  "O5" 1
Set-Num  "P5" 1
Set-Num  "Q5" 1
Set-Text "R5" "Ask"
Set-Link "S5" "artfynd/A 61558-2023 artfynd.xlsx"            "A 61558-2023"
Set-Link "T5" "kartor/A 61558-2023 karta.png"                 "A 61558-2023"
Set-Link "V5" "klagomål/A 61558-2023 FSC-klagomål.docx"       "A 61558-2023"
Set-Link "W5" "klagomålsmail/A 61558-2023 FSC-klagomål mail.docx" "A 61558-2023"
Set-Link "X5" "tillsyn/A 61558-2023 tillsynsbegäran.docx"     "A 61558-2023"
Set-Link "Y5" "tillsynsmail/A 61558-2023 tillsynsbegäran mail.docx" "A 61558-2023"

# --- Row 6 (becomes old row 5: "A 13766-2023") ---
Set-Text "A6" "A 13766-2023"
Set-Num  "B6" 45007
Clear-Cell "F6"
Set-Num  "G6" 0.9
Set-Num  "H6" 1
Set-Num  "I6" 0
Set-Num  "J6" 0
Set-Num  "K6" 0
Set-Num  "L6" 0
Set-Num  "M6" 0
Set-Num  "N6" 0
Set-Num  "O6" 0
Set-Num  "P6" 0
Set-Num  "Q6" 1
Set-Text "R6" "Större vattensalamander"
Set-Link "S6" "artfynd/A 13766-2023 artfynd.xlsx"            "A 13766-2023"
Set-Link "T6" "kartor/A 13766-2023 karta.png"                 "A 13766-2023"
Set-Link "V6" "klagomål/A 13766-2023 FSC-klagomål.docx"       "A 13766-2023"
Set-Link "W6" "klagomålsmail/A 13766-2023 FSC-klagomål mail.docx" "A 13766-2023"
Set-Link "X6" "tillsyn/A 13766-2023 tillsynsbegäran.docx"     "A 13766-2023"
Set-Link "Y6" "tillsynsmail/A 13766-2023 tillsynsbegäran mail.docx" "A 13766-2023"

# --- Row 7 (becomes old row 8: "A 60891-2024") ---
Set-Text "A7" "A 60891-2024"
Set-Num  "B7" 45644
Set-Num  "G7" 16.1
Set-Num  "H7" 1
Set-Num  "I7" 0
Set-Text "R7" "Lövgroda"
Set-Link "S7" "artfynd/A 60891-2024 artfynd.xlsx"            "A 60891-2024"
Set-Link "T7" "kartor/A 60891-2024 karta.png"                 "A 60891-2024"
Set-Link "V7" "klagomål/A 60891-2024 FSC-klagomål.docx"       "A 60891-2024"
Set-Link "W7" "klagomålsmail/A 60891-2024 FSC-klagomål mail.docx" "A 60891-2024"
Set-Link "X7" "tillsyn/A 60891-2024 tillsynsbegäran.docx"     "A 60891-2024"
Set-Link "Y7" "tillsynsmail/A 60891-2024 tillsynsbegäran mail.docx" "A 60891-2024"

# --- Row 8 (becomes old row 7: "A 49546-2025") ---
Set-Text "A8" "A 49546-2025"
Set-Num  "B8" 45939
Set-Num  "G8" 4.4
Set-Num  "H8" 0
Set-Num  "I8" 1
Set-Text "R8" "Igelkottsröksvamp"
Set-Link "S8" "artfynd/A 49546-2025 artfynd.xlsx"            "A 49546-2025"
Set-Link "T8" "kartor/A 49546-2025 karta.png"                 "A 49546-2025"
Set-Link "V8" "klagomål/A 49546-2025 FSC-klagomål.docx"       "A 49546-2025"
Set-Link "W8" "klagomålsmail/A 49546-2025 FSC-klagomål mail.docx" "A 49546-2025"
Set-Link "X8" "tillsyn/A 49546-2025 tillsynsbegäran.docx"     "A 49546-2025"
Set-Link "Y8" "tillsynsmail/A 49546-2025 tillsynsbegäran mail.docx" "A 49546-2025"

# --- Row 14 (becomes old row 15: "A 49536-2025") ---
Set-Text "A14" "A 49536-2025"
Set-Num  "B14" 45939.4221875
Set-Num  "G14" 1.5

# --- Row 15 (becomes old row 16: "A 49543-2025") ---
Set-Text "A15" "A 49543-2025"
Set-Num  "B15" 45939.42862268518
Set-Num  "G15" 1.4

# --- Row 16 (becomes old row 14: "A 32596-2024") ---
Set-Text "A16" "A 32596-2024"
Set-Num  "B16" 45513.61667824074
Set-Num  "G16" 2.6

# --- Row 17 (becomes old row 18: "A 40417-2022") ---
Set-Text "A17" "A 40417-2022"
Set-Num  "B17" 44823
Set-Num  "G17" 2.3

# --- Row 18 (becomes old row 19: "A 18090-2022") ---
Set-Text "A18" "A 18090-2022"
Set-Num  "B18" 44684
Set-Num  "G18" 4.9

# --- Row 19 (becomes old row 20: "A 60803-2023") ---
Set-Text "A19" "A 60803-2023"
Set-Num  "B19" 45260
Set-Num  "G19" 1.6

# --- Row 20 (becomes old row 17: "A 49549-2025") ---
Set-Text "A20" "A 49549-2025"
Set-Num  "B20" 45939
Set-Num  "G20" 0.5
